# Apply the "LOCAL LM" cierre_caja update:
#  - rename the sheet (and its workbook.xml <sheet name=.../> entry) from
#    "Velázquez" to "LOCAL LM"
#  - rewrite the data table (rows 2-9) with new values and blank out column B
#  - append 8 new data rows (10-17)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "LOCAL LM"

# Columns C (dates like "05/10/2024") and G (plain integers like "6") would
# otherwise get auto-converted by Excel into a date serial / a number.
# Force them to stay TEXT so they match the original inline-string cells.
$ws.Range("C2:C17").NumberFormat = "@"
$ws.Range("G2:G17").NumberFormat = "@"

# Tienda, Nombre_TPV, fecha, cierre_tpv_desc, Nombre_MdP, total_ventas, total_operaciones
$rows = @(
    ,@("LOCAL LM", "", "05/10/2024", "Mañana", "GLOVO",        "190,60",  "6")
    ,@("LOCAL LM", "", "05/10/2024", "Mañana", "SMS",          "7,20",    "2")
    ,@("LOCAL LM", "", "05/10/2024", "Mañana", "EUROS",        "1915,04", "182")
    ,@("LOCAL LM", "", "05/10/2024", "Mañana", "TARJETA VISA", "4663,50", "414")
    ,@("LOCAL LM", "", "05/10/2024", "Tarde",  "EUROS",        "295,50",  "32")
    ,@("LOCAL LM", "", "05/10/2024", "Tarde",  "TARJETA VISA", "857,99",  "71")
    ,@("LOCAL LM", "", "05/10/2024", "Mañana", "EUROS",        "563,28",  "46")
    ,@("LOCAL LM", "", "05/10/2024", "Mañana", "GLOVO",        "73,80",   "5")
    ,@("LOCAL LM", "", "05/10/2024", "Mañana", "TARJETA VISA", "1351,34", "90")
    ,@("LOCAL LM", "", "05/10/2024", "Tarde",  "TARJETA VISA", "2196,68", "198")
    ,@("LOCAL LM", "", "05/10/2024", "Tarde",  "EUROS",        "1107,08", "100")
    ,@("LOCAL LM", "", "05/10/2024", "Tarde",  "EUROS",        "563,28",  "46")
    ,@("LOCAL LM", "", "05/10/2024", "Tarde",  "GLOVO",        "73,80",   "5")
    ,@("LOCAL LM", "", "05/10/2024", "Tarde",  "TARJETA VISA", "1351,34", "90")
    ,@("LOCAL LM", "", "05/10/2024", "Mañana", "EUROS",        "295,50",  "32")
    ,@("LOCAL LM", "", "05/10/2024", "Mañana", "TARJETA VISA", "857,99",  "71")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -eq "") {
        $ws.Cells.Item($r, 2).ClearContents()
    } else {
        $ws.Cells.Item($r, 2).Value = $row[1]
    }
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}
